# Update requirements list: replace old "RF" functional requirements with
# new "RNF" (non-functional) requirements, and extend the list with two
# more rows plus a couple of trailing blank (but sized) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("RNF1", "adicionar serviço"),
    @("RNF2", "alterar serviço"),
    @("RNF3", "remover serviço"),
    @("RNF4", "confirmar agendamento"),
    @("RNF5", "criar conta empresa"),
    @("RNF6", "criar conta cliente"),
    @("RNF7", "cadastrar pet"),
    @("RNF9", "atualizar dados do pet"),
    @("RNF10", "remover pet"),
    @("RNF11", "solicitar agendamento")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Rows.Item($row).RowHeight = 21
}

# Row 3 keeps its original (slightly smaller) height.
$ws.Rows.Item(3).RowHeight = 20

# Two trailing blank rows, still explicitly sized like the rest.
$ws.Rows.Item(12).RowHeight = 21
$ws.Rows.Item(13).RowHeight = 21

$ws.Range("A13").Select()
